$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("movies")

# Update the poster_path formula for row 2 (literal, not shared) and
# rows 3-19 (shared formula with anchor on F3).
$ws.Range("F2").Formula = '=CONCATENATE("/static/img/posters/",A2)'
$ws.Range("F3:F19").Formula = '=CONCATENATE("/static/img/posters/",A3)'

# Column F (and its header) moves from a Text-forced format to General.
$ws.Columns.Item(6).NumberFormat = "General"

# Make "movies" the active sheet / tab, landing the selection on G20,
# matching where the user left off after the refactor.
$ws.Activate()
$ws.Range("G20").Select()
